# Applies the "update Sheets via scheduled runner" data refresh:
# recomputed market-price columns (H:N) for a set of Leve rows across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
#
# Sheet order (index -> name): 1=ALC 2=ARM 3=BSM 4=CRP 5=CUL 6=GSM 7=LTW 8=WVR
# Columns: H=8 currentAveragePrice, I=9 currentAveragePriceNQ,
#          J=10 currentAveragePriceHQ, K=11 LevePriceNQ, L=12 LevePriceHQ,
#          M=13 LeveProfitNQ, N=14 LeveProfitHQ
#
# Each data line is: <sheetIndex> <row> <col> <SET|CLEAR> [<newValue>]
# CLEAR means the recomputed cell has no value for this refresh (the
# source cell is removed, matching rows where a profit column is blank).

$wb = $excel.ActiveWorkbook

$cellOps = @"
1	41	8	SET	1637.7778
1	41	9	SET	2696.6667
1	41	10	SET	1108.3334
1	41	11	SET	2696.6667
1	41	12	SET	1108.3334
1	41	13	SET	-2256.6667
1	41	14	SET	-1988.3334
1	57	8	SET	31389.5
1	57	10	SET	31389.5
1	57	12	SET	94168.5
1	57	14	SET	-95166.5
1	58	8	SET	750.3333
1	58	9	SET	250.5
1	58	10	SET	1750
1	58	11	SET	751.5
1	58	12	SET	5250
1	58	13	SET	-601.5
1	58	14	SET	-5550
1	64	8	SET	3300
1	64	9	SET	3100
1	64	11	SET	3100
1	64	13	SET	-2852
1	67	8	SET	3300
1	67	9	SET	3100
1	67	11	SET	3100
1	67	13	SET	-2242
1	107	8	SET	239
1	107	9	SET	267.7
1	107	11	SET	267.7
1	107	13	SET	1652.3
1	135	8	SET	720.2778
1	135	9	SET	691.06665
1	135	11	SET	6219.59985
1	135	13	SET	-3684.59985
1	137	8	SET	48873.19
1	137	9	SET	1003.4167
1	137	11	SET	3010.2501
1	137	13	SET	-460.2501000000002
1	138	8	SET	1767.98
1	138	9	SET	1062.2778
1	138	10	SET	2164.9375
1	138	11	SET	3186.8334
1	138	12	SET	6494.8125
1	138	13	SET	1953.1666
1	138	14	SET	-16774.8125
2	32	8	SET	5556.59
2	32	9	SET	3745.2195
2	32	10	SET	13808.389
2	32	11	SET	3745.2195
2	32	12	SET	13808.389
2	32	13	SET	-3458.2195
2	32	14	SET	-14382.389
2	61	8	SET	40614.855
2	61	9	SET	50950.75
2	61	11	SET	50950.75
2	61	13	SET	-50738.75
2	63	8	SET	2017.8
2	63	9	SET	2017.8
2	63	11	SET	2017.8
2	63	13	SET	-1331.8
2	66	8	SET	2017.8
2	66	9	SET	2017.8
2	66	11	SET	10089
2	66	13	SET	-6657
2	74	8	SET	547.4545000000001
2	74	9	SET	547.4545000000001
2	74	11	SET	547.4545000000001
2	74	13	SET	326.5454999999999
2	77	8	SET	547.4545000000001
2	77	9	SET	547.4545000000001
2	77	11	SET	2737.2725
2	77	13	SET	1630.7275
2	102	8	SET	2518.5
2	102	9	SET	2518.5
2	102	11	SET	2518.5
2	102	13	SET	-896.5
2	136	8	SET	40614.855
2	136	9	SET	50950.75
2	136	11	SET	152852.25
2	136	13	SET	-150302.25
3	86	8	SET	183335
3	86	9	SET	928.6667
3	86	11	SET	928.6667
3	86	13	SET	194.3333
3	89	8	SET	183335
3	89	9	SET	928.6667
3	89	11	SET	4643.3335
3	89	13	SET	972.6665000000003
3	134	8	SET	4190.0894
3	134	9	SET	4373.4634
3	134	11	SET	13120.3902
3	134	13	SET	-10585.3902
4	22	8	SET	1099.8667
4	22	9	SET	441.66666
4	22	10	SET	1538.6666
4	22	11	SET	441.66666
4	22	12	SET	1538.6666
4	22	13	SET	-91.66665999999998
4	22	14	SET	-2238.6666
4	62	8	SET	2662.4
4	62	9	SET	2450
4	62	11	SET	2450
4	62	13	SET	-1826
4	65	8	SET	2662.4
4	65	9	SET	2450
4	65	11	SET	12250
4	65	13	SET	-9130
4	132	8	SET	1800.6111
4	132	9	SET	1211.75
4	132	11	SET	3635.25
4	132	13	SET	-1105.25
4	134	8	SET	1364.92
4	134	9	SET	1193.8823
4	134	11	SET	3581.6469
4	134	13	SET	-1046.6469
5	5	8	SET	569.4167
5	5	9	SET	522.8095
5	5	10	SET	895.6667
5	5	11	SET	1568.4285
5	5	12	SET	2687.0001
5	5	13	SET	-1456.4285
5	5	14	SET	-2911.0001
5	47	8	SET	3429.5
5	47	9	SET	2744.5
5	47	10	SET	4799.5
5	47	11	SET	8233.5
5	47	12	SET	14398.5
5	47	13	SET	-7802.5
5	47	14	SET	-15260.5
5	122	8	SET	1062.3182
5	122	10	SET	1108.7
5	122	12	SET	9978.300000000001
5	122	14	SET	-14878.3
5	131	8	SET	22404.781
5	131	10	SET	29777.291
5	131	12	SET	89331.87300000001
5	131	14	SET	-99411.87300000001
5	135	8	SET	569.4167
5	135	9	SET	522.8095
5	135	10	SET	895.6667
5	135	11	SET	4705.2855
5	135	12	SET	8061.0003
5	135	13	SET	-2170.2855
5	135	14	SET	-13131.0003
5	140	8	SET	4988.5
5	140	9	SET	4520.857
5	140	11	SET	13562.571
5	140	13	SET	-8382.571
6	52	8	SET	55677.668
6	52	10	SET	55677.668
6	52	12	SET	55677.668
6	52	14	SET	-56195.668
6	132	8	SET	2139588.2
6	132	9	SET	2566533.5
6	132	11	SET	7699600.5
6	132	13	SET	-7697070.5
7	69	8	SET	77032.60000000001
7	69	9	SET	0
7	69	10	SET	77032.60000000001
7	69	11	SET	0
7	69	12	SET	77032.60000000001
7	69	13	CLEAR	
7	69	14	SET	-78654.60000000001
7	72	8	SET	77032.60000000001
7	72	9	SET	0
7	72	10	SET	77032.60000000001
7	72	11	SET	0
7	72	12	SET	231097.8
7	72	13	CLEAR	
7	72	14	SET	-239209.8
7	132	8	SET	4052.8484
7	132	9	SET	3120.0625
7	132	10	SET	4930.7646
7	132	11	SET	9360.1875
7	132	12	SET	14792.2938
7	132	13	SET	-6830.1875
7	132	14	SET	-19852.2938
7	136	8	SET	3653.889
7	136	9	SET	2637.4
7	136	10	SET	4924.5
7	136	11	SET	7912.200000000001
7	136	12	SET	14773.5
7	136	13	SET	-5362.200000000001
7	136	14	SET	-19873.5
8	113	8	SET	873.26666
8	113	9	SET	679.4
8	113	11	SET	2038.2
8	113	13	SET	131.8000000000002
8	132	8	SET	1271.4906
8	132	9	SET	1249.2059
8	132	11	SET	3747.6177
8	132	13	SET	-1217.6177
8	136	8	SET	15016296
8	136	9	SET	32680776
8	136	10	SET	1489.05
8	136	11	SET	98042328
8	136	12	SET	4467.15
8	136	13	SET	-98039778
8	136	14	SET	-9567.15
8	137	8	SET	0
8	137	10	SET	0
8	137	12	SET	0
8	137	14	CLEAR	
"@

$lines = $cellOps -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split "`t"
    $sheetIdx = [int]$parts[0]
    $row = [int]$parts[1]
    $col = [int]$parts[2]
    $action = $parts[3]

    $ws = $wb.Worksheets.Item($sheetIdx)
    $cell = $ws.Cells.Item($row, $col)

    if ($action -eq "SET") {
        $cell.Value = [double]$parts[4]
    } elseif ($action -eq "CLEAR") {
        $cell.ClearContents()
    }
}

Write-Host "Applied $($lines.Count) cell updates across 8 sheets"
